$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-02-29 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-01 Friday", 2) | Out-Null

# Update each arithmetic problem cell in the table (scoped per-cell to avoid
# ambiguity since several old/new values repeat across the table)
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("98÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "94÷6=", 2) | Out-Null

$cell = $t.Cell(1, 2)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("37÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "50÷9=", 2) | Out-Null

$cell = $t.Cell(1, 3)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("62÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷8=", 2) | Out-Null

$cell = $t.Cell(1, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("30÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "81÷9=", 2) | Out-Null

$cell = $t.Cell(1, 5)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("88÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷6=", 2) | Out-Null

$cell = $t.Cell(5, 1)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("67÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷9=", 2) | Out-Null

$cell = $t.Cell(5, 2)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("98÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "43÷2=", 2) | Out-Null

$cell = $t.Cell(5, 3)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("86÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "49÷9=", 2) | Out-Null

$cell = $t.Cell(5, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("69÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷9=", 2) | Out-Null

$cell = $t.Cell(5, 5)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("13÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "66÷4=", 2) | Out-Null

$cell = $t.Cell(9, 1)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("86÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "33÷8=", 2) | Out-Null

$cell = $t.Cell(9, 2)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("80÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "53÷3=", 2) | Out-Null

$cell = $t.Cell(9, 3)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("40÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "78÷3=", 2) | Out-Null

$cell = $t.Cell(9, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("26÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷3=", 2) | Out-Null

$cell = $t.Cell(9, 5)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("93÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "25÷2=", 2) | Out-Null

$cell = $t.Cell(13, 1)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("25÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷2=", 2) | Out-Null

$cell = $t.Cell(13, 2)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("60÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷5=", 2) | Out-Null

$cell = $t.Cell(13, 3)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("34÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "77÷2=", 2) | Out-Null

$cell = $t.Cell(13, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("99÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷9=", 2) | Out-Null

$cell = $t.Cell(13, 5)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("65÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷7=", 2) | Out-Null

$cell = $t.Cell(17, 1)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("64÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "60÷8=", 2) | Out-Null

$cell = $t.Cell(17, 2)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("54÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷7=", 2) | Out-Null

$cell = $t.Cell(17, 3)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("24÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "13÷3=", 2) | Out-Null

$cell = $t.Cell(17, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("69÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "43÷6=", 2) | Out-Null

$cell = $t.Cell(17, 5)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Find.Execute("79÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "50÷2=", 2) | Out-Null

